$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price and volume data (GitHub Actions refresh)
# Rows 11/12 and 18/19 also swap position (coin name/link/price/volume)

# Row 2
$ws.Range("D2").Value = "'30.754.52"
$ws.Range("E2").Value = "'  +1.92%  "

# Row 3
$ws.Range("D3").Value = "'1.876.22"
$ws.Range("E3").Value = "'  +2.03%  "

# Row 4
$ws.Range("D4").Value = "'0.9990"

# Row 5
$ws.Range("D5").Value = "'237.34"
$ws.Range("E5").Value = "'  +1.86%  "

# Row 6
$ws.Range("D6").Value = "'0.9995"
$ws.Range("E6").Value = "'  -0.06%  "

# Row 7
$ws.Range("D7").Value = "'0.4788"
$ws.Range("E7").Value = "'  +2.42%  "

# Row 8
$ws.Range("D8").Value = "'0.2837"
$ws.Range("E8").Value = "'  +4.80%  "

# Row 9
$ws.Range("D9").Value = "'0.06532"
$ws.Range("E9").Value = "'  +4.09%  "

# Row 10
$ws.Range("D10").Value = "'18.95"
$ws.Range("E10").Value = "'  +18.02%  "

# Row 11
$ws.Range("B11").Value = "WrappedEther"
$ws.Range("C11").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D11").Value = "'1.882.08"
$ws.Range("E11").Value = "'  +2.42%  "

# Row 12
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").Value = "'0.07515"
$ws.Range("E12").Value = "'  +1.38%  "

# Row 13
$ws.Range("D13").Value = "'94.78"
$ws.Range("E13").Value = "'  +13.19%  "

# Row 14
$ws.Range("D14").Value = "'5.100"
$ws.Range("E14").Value = "'  +3.45%  "

# Row 15
$ws.Range("D15").Value = "'0.6524"
$ws.Range("E15").Value = "'  +5.16%  "

# Row 16
$ws.Range("D16").Value = "'296.70"
$ws.Range("E16").Value = "'  +30.10%  "

# Row 17
$ws.Range("D17").Value = "'30.703.08"
$ws.Range("E17").Value = "'  +2.05%  "

# Row 18
$ws.Range("B18").Value = "Avalanche"
$ws.Range("C18").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D18").Value = "'13.09"
$ws.Range("E18").Value = "'  +6.09%  "

# Row 19
$ws.Range("B19").Value = "Dai"
$ws.Range("C19").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D19").Value = "'0.9996"
$ws.Range("E19").Value = "'  -0.04%  "

# Row 20
$ws.Range("D20").Value = "'0.000007544"
$ws.Range("E20").Value = "'  +3.68%  "

# Row 21
$ws.Range("D21").Value = "'2.107.40"
$ws.Range("E21").Value = "'  +1.41%  "

# Row 22
$ws.Range("D22").Value = "'0.9989"
$ws.Range("E22").Value = "'  -0.08%  "

# Row 23
$ws.Range("D23").Value = "'5.193"
$ws.Range("E23").Value = "'  +6.63%  "

# Row 24
$ws.Range("D24").Value = "'6.106"
$ws.Range("E24").Value = "'  +4.53%  "

# Row 25
$ws.Range("D25").Value = "'169.27"
$ws.Range("E25").Value = "'  +3.15%  "

# Row 26
$ws.Range("D26").Value = "'9.255"
$ws.Range("E26").Value = "'  +0.64%  "

# Row 27
$ws.Range("D27").Value = "'19.78"
$ws.Range("E27").Value = "'  +11.04%  "

# Row 28
$ws.Range("D28").Value = "'1.984"
$ws.Range("E28").Value = "'  +5.32%  "

# Row 29
$ws.Range("D29").Value = "'0.1054"
$ws.Range("E29").Value = "'  +2.51%  "

# Row 30
$ws.Range("D30").Value = "'1.352"
$ws.Range("E30").Value = "'  -1.38%  "

# Row 31
$ws.Range("D31").Value = "'4.140"
$ws.Range("E31").Value = "'  +1.69%  "

# Row 32
$ws.Range("D32").Value = "'3.960"
$ws.Range("E32").Value = "'  +4.48%  "

# Row 33
$ws.Range("D33").Value = "'0.04986"
$ws.Range("E33").Value = "'  +3.82%  "

# Row 34
$ws.Range("D34").Value = "'1.179"
$ws.Range("E34").Value = "'  +3.60%  "

# Row 35
$ws.Range("D35").Value = "'0.7250"
$ws.Range("E35").Value = "'  +2.41%  "

# Row 36
$ws.Range("E36").Value = "'  +0.47%  "

# Row 37
$ws.Range("D37").Value = "'0.01932"
$ws.Range("E37").Value = "'  +2.54%  "

# Row 38
$ws.Range("D38").Value = "'2.701"
$ws.Range("E38").Value = "'  +1.89%  "

# Row 39
$ws.Range("D39").Value = "'2.056"
$ws.Range("E39").Value = "'  +7.07%  "

# Row 40
$ws.Range("D40").Value = "'0.8917"
$ws.Range("E40").Value = "'  -0.29%  "

# Row 41
$ws.Range("D41").Value = "'107.66"
$ws.Range("E41").Value = "'  +3.21%  "

# Row 42
$ws.Range("D42").Value = "'0.9998"
$ws.Range("E42").Value = "'  -0.18%  "

# Row 43
$ws.Range("D43").Value = "'0.4217"
$ws.Range("E43").Value = "'  +5.43%  "

# Row 44
$ws.Range("D44").Value = "'5.579"
$ws.Range("E44").Value = "'  +0.85%  "

# Row 45
$ws.Range("D45").Value = "'7.374"
$ws.Range("E45").Value = "'  +5.29%  "

# Row 46
$ws.Range("D46").Value = "'65.69"
$ws.Range("E46").Value = "'  +10.21%  "

# Row 47
$ws.Range("D47").Value = "'0.1231"
$ws.Range("E47").Value = "'  +3.32%  "

# Row 48
$ws.Range("D48").Value = "'34.68"
$ws.Range("E48").Value = "'  +6.09%  "

# Row 49
$ws.Range("D49").Value = "'8.787"
$ws.Range("E49").Value = "'  +2.85%  "

# Row 50
$ws.Range("D50").Value = "'1.397"
$ws.Range("E50").Value = "'  +2.94%  "

# Row 51
$ws.Range("D51").Value = "'0.05554"
$ws.Range("E51").Value = "'  +0.89%  "

# Clear the auto-applied text-format styling from the apostrophe-prefixed assignments above
$ws.Range("D2:E51").ClearFormats()
